# Insert a new weekly data row at row 117 (pushing existing rows 117-187
# down to 118-188) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(117).Insert()

$ws.Cells.Item(117, 1).Value = 7
$ws.Cells.Item(117, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(117, 3).Value = "Ñuble"
$ws.Cells.Item(117, 4).Value = 44488
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 100114013
$ws.Cells.Item(117, 7).Value = "Zanahoria"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 100
$ws.Cells.Item(117, 11).Value = 8500
$ws.Cells.Item(117, 12).Value = 9000
$ws.Cells.Item(117, 13).Value = 8750
$ws.Cells.Item(117, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(117, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(117, 16).Value = 438
$ws.Cells.Item(117, 17).Value = 20
$ws.Cells.Item(117, 18).Value = "Hortaliza"
